$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.5202807690922467
$ws.Range("J2").Value = 0.5202807690922467
$ws.Range("M2").Value = 0.602545
$ws.Range("N2").Value = 1.807635
$ws.Range("O2").Value = 0.1012915997907189
$ws.Range("P2").Value = 0.1012915997907189
$ws.Range("Q2").Value = 0.243720414325
$ws.Range("R2").Value = 2.193483728925
$ws.Range("S2").Value = 0.05270007144169928
$ws.Range("T2").Value = 0.05270007144169928
$ws.Range("I3").Value = 0.5202807690922467
$ws.Range("J3").Value = 0.5202807690922467
$ws.Range("O3").Value = 0.402828040777877
$ws.Range("P3").Value = 0.4028280407778771
$ws.Range("S3").Value = 0.2095836828678368
$ws.Range("T3").Value = 0.2095836828678369
$ws.Range("I4").Value = 0.5202807690922467
$ws.Range("J4").Value = 0.5202807690922467
$ws.Range("M4").Value = 2.949802666666667
$ws.Range("N4").Value = 8.849408
$ws.Range("O4").Value = 0.495880359431404
$ws.Range("P4").Value = 0.4958803594314041
$ws.Range("Q4").Value = 1.193150931626667
$ws.Range("R4").Value = 10.73835838464
$ws.Range("S4").Value = 0.2579970147827106
$ws.Range("T4").Value = 0.2579970147827107
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.372951
$ws.Range("H5").Value = 1.118853
$ws.Range("I5").Value = 0.4797192309077532
$ws.Range("J5").Value = 0.4797192309077531
$ws.Range("M5").Value = 0.602545
$ws.Range("N5").Value = 1.807635
$ws.Range("O5").Value = 0.1012915997907189
$ws.Range("P5").Value = 0.1012915997907189
$ws.Range("Q5").Value = 0.224719760295
$ws.Range("R5").Value = 2.022477842655
$ws.Range("S5").Value = 0.04859152834901959
$ws.Range("T5").Value = 0.04859152834901959
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.372951
$ws.Range("H6").Value = 1.118853
$ws.Range("I6").Value = 0.4797192309077532
$ws.Range("J6").Value = 0.4797192309077531
$ws.Range("O6").Value = 0.402828040777877
$ws.Range("P6").Value = 0.4028280407778771
$ws.Range("Q6").Value = 0.89369129277
$ws.Range("R6").Value = 8.043221634930001
$ws.Range("S6").Value = 0.1932443579100402
$ws.Range("T6").Value = 0.1932443579100402
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.372951
$ws.Range("H7").Value = 1.118853
$ws.Range("I7").Value = 0.4797192309077532
$ws.Range("J7").Value = 0.4797192309077531
$ws.Range("M7").Value = 2.949802666666667
$ws.Range("N7").Value = 8.849408
$ws.Range("O7").Value = 0.495880359431404
$ws.Range("P7").Value = 0.4958803594314041
$ws.Range("Q7").Value = 1.100131854336
$ws.Range("R7").Value = 9.901186689024001
$ws.Range("S7").Value = 0.2378833446486933
$ws.Range("T7").Value = 0.2378833446486934